# Insert a new "DormitoryAmount" column between DormitoryAddress (B) and
# RoomNumber (C), shifting the old C/D (RoomNumber/Capacity) columns to D/E.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift existing RoomNumber/Capacity columns (C:D) one column to the right
# so there is room for the new DormitoryAmount column at C.
$ws.Columns("C:C").Insert() | Out-Null

# Header for the newly inserted column.
$ws.Range("C1").Value = "DormitoryAmount"

# Fill in the DormitoryAmount values for each data row.
$ws.Range("C2").Value = 200
$ws.Range("C3").Value = 200
$ws.Range("C4").Value = 200
$ws.Range("C5").Value = 300
$ws.Range("C6").Value = 300
